$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.1286665469352575
$ws.Range("C2").Value = 2.584314119446431
$ws.Range("D2").Value = 11.86673044746297
$ws.Range("E2").Value = 3.444812106263994
$ws.Range("F2").Value = 3.475993634475042
$ws.Range("G2").Value = 52

$ws.Range("B3").Value = 0.208472532147568
$ws.Range("C3").Value = 1.905008350356748
$ws.Range("D3").Value = 5.6589047777894
$ws.Range("E3").Value = 2.378845261421894
$ws.Range("F3").Value = 2.393272416744899
$ws.Range("G3").Value = 51

$ws.Range("B4").Value = 0.3820242035387909
$ws.Range("C4").Value = 1.902010132220234
$ws.Range("D4").Value = 5.867547615856401
$ws.Range("E4").Value = 2.422302131414742
$ws.Range("F4").Value = 2.416272454685605
$ws.Range("G4").Value = 50

$ws.Range("B5").Value = 0.191239345277837
$ws.Range("C5").Value = 1.716821400939533
$ws.Range("D5").Value = 5.499592484960565
$ws.Range("E5").Value = 2.345120995803962
$ws.Range("F5").Value = 2.36634644499125
$ws.Range("G5").Value = 41

$ws.Range("B6").Value = 0.2671928775352811
$ws.Range("C6").Value = 1.741437505328719
$ws.Range("D6").Value = 5.032127257509086
$ws.Range("E6").Value = 2.243240347691055
$ws.Range("F6").Value = 2.264087689223547
$ws.Range("G6").Value = 31

$ws.Range("B7").Value = 0.3126204252219001
$ws.Range("C7").Value = 1.694986524670226
$ws.Range("D7").Value = 4.73978801325031
$ws.Range("E7").Value = 2.177105420793929
$ws.Range("F7").Value = 2.192679688080215
$ws.Range("G7").Value = 29

$ws.Range("B8").Value = 0.3719975688557745
$ws.Range("C8").Value = 1.867752137679124
$ws.Range("D8").Value = 5.380359802164317
$ws.Range("E8").Value = 2.319560260515841
$ws.Range("F8").Value = 2.333150688324054
$ws.Range("G8").Value = 27

$ws.Range("B9").Value = 0.4130801636093303
$ws.Range("C9").Value = 2.085153282719833
$ws.Range("D9").Value = 6.310345893095409
$ws.Range("E9").Value = 2.512040185406159
$ws.Range("F9").Value = 2.545742663513931
$ws.Range("G9").Value = 19

$ws.Range("B10").Value = 1.987269086558114
$ws.Range("C10").Value = 1.995519100223671
$ws.Range("D10").Value = 6.088527956227622
$ws.Range("E10").Value = 2.46749426670613
$ws.Range("F10").Value = 1.527668288782175
$ws.Range("G10").Value = 12

$ws.Range("B11").Value = 1.658904638620646
$ws.Range("C11").Value = 1.679380840373142
$ws.Range("D11").Value = 3.778169217003525
$ws.Range("E11").Value = 1.943751325916853
$ws.Range("F11").Value = 1.132588085407946
$ws.Range("G11").Value = 5
